$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 (MAG file name)
$ws.Range("A2").Value = "even_MAG-GUT33657.fa"

# Update numeric prediction values for row 2 (B2:T2)
$ws.Range("B2").Value = 0.01416995289984451
$ws.Range("C2").Value = 0.007032596430465727
$ws.Range("D2").Value = 0.001914900352217104
$ws.Range("E2").Value = 0.02636301097311251
$ws.Range("F2").Value = 0.001769329767789052
$ws.Range("G2").Value = 0.01698750834484411
$ws.Range("H2").Value = 0.006737378631208064
$ws.Range("I2").Value = 0.004335880739152437
$ws.Range("J2").Value = 0.002012529197051532
$ws.Range("K2").Value = 0.01324921019105777
$ws.Range("L2").Value = [double]"4.664237198882104e-05"
$ws.Range("M2").Value = [double]"6.960040999231455e-05"
$ws.Range("N2").Value = 0.004325296250001704
$ws.Range("O2").Value = 0.3857400106137104
$ws.Range("P2").Value = 0.5131914072144675
$ws.Range("Q2").Value = 0.001987479904353319
$ws.Range("R2").Value = [double]"6.726569376910275e-05"
$ws.Range("S2").Value = [double]"1.497394356299093e-11"
$ws.Range("T2").Value = 0.5131914072144675

# Add new column V: header "rejection-f" with same formatting as U1 (copy style), plus row 2 value
$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("V1").Value = "rejection-f"

$ws.Range("V2").Value = "s__Enterocloster sp900549235(reject)"
